# Job parser and tweaks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shee1")
$ws.Activate()

# Rename header row to camelCase job-parser field names
$ws.Range("A1").Value = "brand"
$ws.Range("B1").Value = "position"
$ws.Range("C1").Value = "address"
$ws.Range("D1").Value = "experience"
$ws.Range("E1").Value = "skills"
$ws.Range("F1").Value = "availability"
$ws.Range("G1").Value = "startTime"
$ws.Range("H1").Value = "compensation"
$ws.Range("I1").Value = "description"
$ws.Range("J1").Value = "openCall"
$ws.Range("K1").Value = "interviewQuestion"
$ws.Range("L1").Value = "trackingUrl"
$ws.Range("M1").Value = "schedule"

# Column M width (COM ColumnWidth uses character units with pixel rounding;
# 11.17 round-trips to a stored OOXML width of 12)
$ws.Columns.Item(13).ColumnWidth = 11.17

# Sample data row
$ws.Range("A2").Value = 123123
$ws.Range("B2").Value = 131
$ws.Range("C2").Value = 132312
$ws.Range("D2").Value = 31
$ws.Range("E2").Value = "skill1|skill2|skill3|skill4"
$ws.Range("F2").Value = "part|full|anytime"

# Move the selection
$ws.Range("F3").Select()
